# v1.4 - Set reviewer verification to closed
# LH_TC_LOGIN_Review_004
#
# 1. Add a new "v1.4" row to the Version History table (row 6, which was
#    already a blank row inside the Table1 range A1:D9).
# 2. Flip "Reviewer verification" (column J) from Open to Closed for the
#    LH_REVIEW_TC_LOGIN_014..019 rows (rows 15-20) on the Review sheet.

$wb = $excel.ActiveWorkbook

# ---- Version History sheet -------------------------------------------------
$history = $wb.Worksheets.Item("Version History")

$history.Range("A6").Value2 = "v1.4"
$history.Range("B6").Value2 = "Mahmoud abdelmageed"
$history.Range("C6").Value2 = "Set Reviewer verification to closed"
$history.Range("D6").Value2 = 45770

# ---- Review sheet -----------------------------------------------------------
$review = $wb.Worksheets.Item("Review sheet")

$review.Range("J15").Value2 = "Closed"
$review.Range("J16").Value2 = "Closed"
$review.Range("J17").Value2 = "Closed"
$review.Range("J18").Value2 = "Closed"
$review.Range("J19").Value2 = "Closed"
$review.Range("J20").Value2 = "Closed"
